$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 1: remove the old summary numbers ("total"/"free" labels, the 2.44% figure) -
# keep F1 ("cb6544b") untouched.
$ws.Range("B1:D1").ClearContents()

# Row 2: the 168792 / 34456 counters are removed entirely.
$ws.Range("C2:D2").ClearContents()

# Row 3: drop the "nothing" column header, keep word/correct/wrong.
$ws.Range("E3").ClearContents()

# Row 42: the stray "x" mark was in column E; move it one column left to D,
# matching the rest of the table's layout.
$ws.Range("D42").Value2 = $ws.Range("E42").Value2
$ws.Range("E42").ClearContents()

# Reflect the author's scroll position / active selection in the saved view.
$ws.Range("E9").Select()
